# Label BOM items better.
# Relabel a handful of package/description strings in the BOM sheet so the
# shared-string text matches the clearer naming scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "R-W4" -> "R-1/4W" for every resistor row that used the old package label
$ws.Range("C2:C9").Value = "R-1/4W"

# Ceramic capacitor package/description
$ws.Range("C12").Value = "C-P5mm"
$ws.Range("E12").Value = "Ceramic Capacitor THT"

# Electrolytic (polarized) capacitor package/description
$ws.Range("C13").Value = "E-P2.5mm 6.3x11.5mm"
$ws.Range("E13").Value = "Electrolytic Capacitor THT"

# Move the active selection to match the saved view state
[void]$ws.Range("E12").Select()
